$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 45183
$ws.Range("C3").Value = 45183
$ws.Range("C4").Value = 45183
$ws.Range("C5").Value = 45183
